$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '63.063.07'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -0.64%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.550.98'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +0.36%  '

$ws.Range("E4").Value = '  -0.06%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '581.45'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +2.43%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '147.10'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -2.16%  '

$ws.Range("E7").Value = '  -0.05%  '

$ws.Range("E8").Value = '  -0.16%  '

$ws.Range("E9").Value = '  +0.31%  '

$ws.Range("E10").Value = '  -2.08%  '

$ws.Range("E11").Value = '  +0.01%  '

$ws.Range("E12").Value = '  -1.07%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '27.53'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -2.31%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '3.006.42'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +0.29%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '62.970.24'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -0.68%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.0000145'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +1.26%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '2.536.27'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -0.13%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '11.33'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -2.10%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '338.73'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -0.04%  '

$ws.Range("E20").Value = '  -0.28%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.74'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -0.63%  '

$ws.Range("E22").Value = '  -0.04%  '

$ws.Range("E23").Value = '  -0.64%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.677.50'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +0.32%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.170'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +0.37%  '

$ws.Range("E26").Value = '  +1.41%  '

$ws.Range("B27").Value = 'Binance-PegBSC-USD'
$ws.Range("C27").Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.999'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -0.21%  '

$ws.Range("B28").Value = 'SuiNetwork'
$ws.Range("C28").Value = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.48'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -4.45%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '8.34'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -1.09%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '7.72'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +8.18%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.94'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +4.34%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.0₃0817'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +0.08%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '178.36'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +0.34%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '422.30'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +0.65%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.56'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -1.30%  '

$ws.Range("E36").Value = '  -1.39%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '19.08'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +0.36%  '

$ws.Range("E39").Value = '  -1.00%  '

$ws.Range("E40").Value = '  -1.06%  '

$ws.Range("E41").Value = '  +0.02%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '39.73'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +0.61%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '150.54'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -1.89%  '

$ws.Range("E44").Value = '  -0.12%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '20.78'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +0.30%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.0539'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +2.75%  '

$ws.Range("E48").Value = '  +0.54%  '

$ws.Range("E49").Value = '  +0.60%  '

$ws.Range("E50").Value = '  -1.77%  '

$ws.Range("E51").Value = '  -5.50%  '
